# Update the "Target vs Sales" branch-wise stock status chart source data
# (rows 2-32, columns B,C,D,E,F,G,H) with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new values for columns B,C,D,E,F,G,H
$newData = @{
    2  = @(76, 0, 24, 0, 0, 0, 52)
    3  = @(76, 0, 31, 10, 6, 3, 26)
    4  = @(76, 0, 35, 14, 4, 3, 20)
    5  = @(76, 0, 34, 7, 7, 7, 21)
    6  = @(76, 0, 29, 10, 3, 10, 24)
    7  = @(76, 0, 30, 11, 3, 6, 26)
    8  = @(76, 0, 32, 13, 4, 4, 23)
    9  = @(76, 0, 30, 12, 2, 5, 27)
    10 = @(76, 0, 30, 5, 6, 3, 32)
    11 = @(76, 0, 32, 17, 1, 7, 19)
    12 = @(76, 0, 29, 14, 5, 7, 21)
    13 = @(76, 0, 33, 9, 4, 5, 25)
    14 = @(76, 0, 35, 10, 6, 6, 19)
    15 = @(76, 0, 37, 6, 4, 4, 25)
    16 = @(76, 0, 31, 8, 4, 3, 30)
    17 = @(76, 0, 27, 7, 5, 8, 29)
    18 = @(76, 0, 30, 10, 2, 10, 24)
    19 = @(76, 0, 35, 11, 3, 4, 23)
    20 = @(76, 0, 36, 7, 1, 4, 28)
    21 = @(76, 0, 33, 8, 7, 6, 22)
    22 = @(76, 0, 34, 14, 1, 6, 21)
    23 = @(76, 0, 33, 9, 3, 7, 24)
    24 = @(76, 0, 32, 19, 6, 6, 13)
    25 = @(76, 0, 34, 8, 9, 5, 20)
    26 = @(76, 0, 31, 7, 2, 7, 29)
    27 = @(76, 0, 27, 9, 3, 5, 32)
    28 = @(76, 0, 28, 11, 2, 6, 29)
    29 = @(76, 0, 29, 8, 2, 8, 29)
    30 = @(76, 0, 30, 12, 4, 8, 22)
    31 = @(76, 0, 34, 8, 5, 9, 20)
    32 = @(76, 0, 31, 12, 2, 8, 23)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    # Columns B through H are columns 2 through 8
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
